# Applies numeric corrections to the Leve profit-tracking tables across all crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1824.6666
$ws.Range("J17").Value = 1989.8
$ws.Range("L17").Value = 5969.4
$ws.Range("N17").Value = -6305.4

$ws.Range("H28").Value = 614.75
$ws.Range("I28").Value = 719.6667
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 719.6667
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = -234.6667
$ws.Range("N28").Value = -1270

$ws.Range("H32").Value = 4999
$ws.Range("I32").Value = 4998.5
$ws.Range("K32").Value = 4998.5
$ws.Range("M32").Value = -4672.5

$ws.Range("H40").Value = 2206.5
$ws.Range("I40").Value = 1936.875
$ws.Range("J40").Value = 2566
$ws.Range("K40").Value = 1936.875
$ws.Range("L40").Value = 2566
$ws.Range("M40").Value = -1761.875
$ws.Range("N40").Value = -2916

$ws.Range("H43").Value = 4832.5557
$ws.Range("I43").Value = 1874.5
$ws.Range("K43").Value = 1874.5
$ws.Range("M43").Value = -1805.5

$ws.Range("H51").Value = 3249.25
$ws.Range("J51").Value = 2999.6667
$ws.Range("L51").Value = 2999.6667
$ws.Range("N51").Value = -3967.6667

$ws.Range("H74").Value = 4950.375
$ws.Range("I74").Value = 2800.4285
$ws.Range("K74").Value = 2800.4285
$ws.Range("M74").Value = -1864.4285

$ws.Range("H77").Value = 4950.375
$ws.Range("I77").Value = 2800.4285
$ws.Range("K77").Value = 14002.1425
$ws.Range("M77").Value = -9322.1425

$ws.Range("H87").Value = 79999
$ws.Range("J87").Value = 79999
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495

$ws.Range("H90").Value = 79999
$ws.Range("J90").Value = 79999
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477

$ws.Range("H99").Value = 191.4
$ws.Range("I99").Value = 191.4
$ws.Range("K99").Value = 574.2
$ws.Range("M99").Value = 923.8

$ws.Range("H116").Value = 7747.5
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H132").Value = 1431.8
$ws.Range("I132").Value = 1347.619
$ws.Range("J132").Value = 1873.75
$ws.Range("K132").Value = 4042.857
$ws.Range("L132").Value = 5621.25
$ws.Range("M132").Value = -1512.857
$ws.Range("N132").Value = -10681.25

$ws.Range("H137").Value = 3027.6428
$ws.Range("I137").Value = 1550.3334
$ws.Range("J137").Value = 4135.625
$ws.Range("K137").Value = 4651.0002
$ws.Range("L137").Value = 12406.875
$ws.Range("M137").Value = -2101.0002
$ws.Range("N137").Value = -17506.875

$ws.Range("H138").Value = 5282.8867
$ws.Range("I138").Value = 2064.6667
$ws.Range("J138").Value = 5475.98
$ws.Range("K138").Value = 6194.000100000001
$ws.Range("L138").Value = 16427.94
$ws.Range("M138").Value = -1054.000100000001
$ws.Range("N138").Value = -26707.94

$ws.Range("H141").Value = 4700
$ws.Range("J141").Value = 7000
$ws.Range("L141").Value = 21000
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1957.7368
$ws.Range("I61").Value = 1585.5
$ws.Range("K61").Value = 1585.5
$ws.Range("M61").Value = -1373.5

$ws.Range("H112").Value = 99999
$ws.Range("J112").Value = 99999
$ws.Range("L112").Value = 99999
$ws.Range("N112").Value = -102953

$ws.Range("H136").Value = 1957.7368
$ws.Range("I136").Value = 1585.5
$ws.Range("K136").Value = 4756.5
$ws.Range("M136").Value = -2206.5

$ws.Range("H139").Value = 89500
$ws.Range("J139").Value = 89500
$ws.Range("L139").Value = 89500
$ws.Range("N139").Value = -99780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 996.5
$ws.Range("I99").Value = 996.5
$ws.Range("K99").Value = 996.5
$ws.Range("M99").Value = 501.5

$ws.Range("H105").Value = 3608.162
$ws.Range("I105").Value = 2846.2693
$ws.Range("K105").Value = 2846.2693
$ws.Range("M105").Value = -1099.2693

$ws.Range("H107").Value = 1248.5
$ws.Range("J107").Value = 1497.5
$ws.Range("L107").Value = 1497.5
$ws.Range("N107").Value = -5337.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 48.307693
$ws.Range("J7").Value = 15.333333
$ws.Range("L7").Value = 15.333333
$ws.Range("N7").Value = -241.333333

$ws.Range("H35").Value = 583.3333
$ws.Range("I35").Value = 583.3333
$ws.Range("K35").Value = 583.3333
$ws.Range("M35").Value = -289.3333

$ws.Range("H99").Value = 10871.96
$ws.Range("I99").Value = 6930.1763
$ws.Range("K99").Value = 6930.1763
$ws.Range("M99").Value = -5432.1763

$ws.Range("H122").Value = 7590.6665
$ws.Range("I122").Value = 6881.375
$ws.Range("J122").Value = 9009.25
$ws.Range("K122").Value = 20644.125
$ws.Range("L122").Value = 27027.75
$ws.Range("M122").Value = -18194.125
$ws.Range("N122").Value = -31927.75

$ws.Range("H126").Value = 10871.96
$ws.Range("I126").Value = 6930.1763
$ws.Range("K126").Value = 20790.5289
$ws.Range("M126").Value = -18320.5289

$ws.Range("H132").Value = 3809.389
$ws.Range("I132").Value = 3078.0833
$ws.Range("K132").Value = 9234.249899999999
$ws.Range("M132").Value = -6704.249899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1751.7273
$ws.Range("I14").Value = 1751.7273
$ws.Range("K14").Value = 5255.1819
$ws.Range("M14").Value = -5082.1819

$ws.Range("H34").Value = 1535.762
$ws.Range("I34").Value = 1139.125
$ws.Range("J34").Value = 1779.8462
$ws.Range("K34").Value = 3417.375
$ws.Range("L34").Value = 5339.5386
$ws.Range("M34").Value = -3333.375
$ws.Range("N34").Value = -5507.5386

$ws.Range("H81").Value = 1921.1666
$ws.Range("J81").Value = 1923.4
$ws.Range("L81").Value = 5770.200000000001
$ws.Range("N81").Value = -8016.200000000001

$ws.Range("H84").Value = 1921.1666
$ws.Range("J84").Value = 1923.4
$ws.Range("L84").Value = 17310.6
$ws.Range("N84").Value = -28542.6

$ws.Range("H113").Value = 902.97144
$ws.Range("J113").Value = 786.3929000000001
$ws.Range("L113").Value = 2359.1787
$ws.Range("N113").Value = -6699.1787

$ws.Range("H129").Value = 7304.4546
$ws.Range("J129").Value = 9237.5
$ws.Range("L129").Value = 27712.5
$ws.Range("N129").Value = -37712.5

$ws.Range("H137").Value = 3446.875
$ws.Range("J137").Value = 4125.6665
$ws.Range("L137").Value = 12376.9995
$ws.Range("N137").Value = -22576.9995

$ws.Range("H138").Value = 4703
$ws.Range("J138").Value = 8266.200000000001
$ws.Range("L138").Value = 24798.6
$ws.Range("N138").Value = -35078.60000000001

$ws.Range("H141").Value = 10507.777
$ws.Range("I141").Value = 10507.777
$ws.Range("K141").Value = 31523.331
$ws.Range("M141").Value = -26343.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 628.5714
$ws.Range("I22").Value = 628.5714
$ws.Range("K22").Value = 628.5714
$ws.Range("M22").Value = -333.5714

$ws.Range("H27").Value = 628.5714
$ws.Range("I27").Value = 628.5714
$ws.Range("K27").Value = 628.5714
$ws.Range("M27").Value = -521.5714

$ws.Range("H55").Value = 379.94736
$ws.Range("I55").Value = 331.93332
$ws.Range("K55").Value = 331.93332
$ws.Range("M55").Value = -158.93332

$ws.Range("H61").Value = 3543.5557
$ws.Range("I61").Value = 3385.4666
$ws.Range("J61").Value = 4334
$ws.Range("K61").Value = 3385.4666
$ws.Range("L61").Value = 4334
$ws.Range("M61").Value = -3183.4666
$ws.Range("N61").Value = -4738

$ws.Range("H100").Value = 4100
$ws.Range("I100").Value = 2010
$ws.Range("K100").Value = 2010
$ws.Range("M100").Value = -1469

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H113").Value = 3543.5557
$ws.Range("I113").Value = 3385.4666
$ws.Range("J113").Value = 4334
$ws.Range("K113").Value = 3385.4666
$ws.Range("L113").Value = 4334
$ws.Range("M113").Value = -1215.4666
$ws.Range("N113").Value = -8674

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6092.6665
$ws.Range("I62").Value = 1998
$ws.Range("J62").Value = 8140
$ws.Range("K62").Value = 1998
$ws.Range("L62").Value = 8140
$ws.Range("M62").Value = -1374
$ws.Range("N62").Value = -9388

$ws.Range("H65").Value = 6092.6665
$ws.Range("I65").Value = 1998
$ws.Range("J65").Value = 8140
$ws.Range("K65").Value = 9990
$ws.Range("L65").Value = 40700
$ws.Range("M65").Value = -6870
$ws.Range("N65").Value = -46940

$ws.Range("H113").Value = 1509.3334
$ws.Range("I113").Value = 1378.1177
$ws.Range("J113").Value = 1648.75
$ws.Range("K113").Value = 4134.3531
$ws.Range("L113").Value = 4946.25
$ws.Range("M113").Value = -1964.3531
$ws.Range("N113").Value = -9286.25

$ws.Range("H122").Value = 1900.7142
$ws.Range("I122").Value = 1917.5
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5752.5
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -3302.5
$ws.Range("N122").Value = -10300

$ws.Range("H136").Value = 45112.434
$ws.Range("J136").Value = 127206.5
$ws.Range("L136").Value = 381619.5
$ws.Range("N136").Value = -386719.5
